# Updates cryptocurrency price (column D) and 1h volume change (column E)
# values on Sheet1, matching the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '64.777.76'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  +1.03%  '

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '3.156.75'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("E4").Value = '  -0.04%  '

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '573.08'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  +2.17%  '

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '151.11'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  +4.71%  '

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.999'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '3.154.72'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  +2.13%  '

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.528'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  +4.55%  '

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '0.163'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  +6.24%  '

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '6.19'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  +0.62%  '

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '0.504'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  +7.21%  '

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0000257'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +12.70%  '

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '37.93'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  +8.12%  '

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '3.667.65'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  +2.04%  '

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '64.865.18'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  +1.08%  '

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '7.24'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  +7.13%  '

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '3.160.30'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("E19").Value = '  +0.46%  '

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '513.14'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  +6.49%  '

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '15.00'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  +7.44%  '

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '0.737'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  +9.51%  '

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '15.14'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +6.90%  '

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '7.86'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  +3.78%  '

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '85.26'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  +5.07%  '

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '0.999'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("E27").Value = '  +4.63%  '

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '8.74'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  +9.65%  '

$ws.Range("E29").Value = '  +5.95%  '

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '28.05'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  +6.58%  '

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '0.999'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("E32").Value = '  +3.40%  '

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '2.66'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  +7.14%  '

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '6.12'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  +9.70%  '

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '6.58'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  +6.06%  '

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '55.69'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  -0.22%  '

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '485.67'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  +7.10%  '

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0424'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  +3.76%  '

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '3.03'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  +2.11%  '

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '3.115.22'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  +4.56%  '

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '8.66'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  +4.97%  '

$ws.Range("E43").Value = '  +5.09%  '

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '0.295'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  +13.32%  '

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '2.47'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  +15.84%  '

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '29.23'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  +4.70%  '

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₃0582'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  +13.15%  '

$ws.Range("E49").Value = '  +3.55%  '

$ws.Range("E50").Value = '  +10.17%  '

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '118.82'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -0.92%  '

